$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I13").Value = "sv"
$ws.Range("J13").Value = "Statement-opinion"
$ws.Range("I32").Value = "sv"
$ws.Range("J32").Value = "Statement-opinion"
$ws.Range("I34").Value = "sv"
$ws.Range("J34").Value = "Statement-opinion"
$ws.Range("I43").Value = "aa"
$ws.Range("J43").Value = "Agree/Accept"
$ws.Range("I45").Value = "b"
$ws.Range("J45").Value = "Acknowledge (Backchannel)"
$ws.Range("I57").Value = "%"
$ws.Range("J57").Value = "Uninterpretable"
$ws.Range("I58").Value = "sd"
$ws.Range("J58").Value = "Statement-non-opinion"
$ws.Range("I74").Value = "%"
$ws.Range("J74").Value = "Uninterpretable"
$ws.Range("I90").Value = "b"
$ws.Range("J90").Value = "Acknowledge (Backchannel)"
$ws.Range("I95").Value = "b"
$ws.Range("J95").Value = "Acknowledge (Backchannel)"
$ws.Range("I117").Value = "sv"
$ws.Range("J117").Value = "Statement-opinion"
$ws.Range("I131").Value = "ba"
$ws.Range("J131").Value = "Appreciation"
$ws.Range("I147").Value = "sd"
$ws.Range("J147").Value = "Statement-non-opinion"
$ws.Range("I153").Value = "b"
$ws.Range("J153").Value = "Acknowledge (Backchannel)"
$ws.Range("I163").Value = "sd"
$ws.Range("J163").Value = "Statement-non-opinion"
$ws.Range("I169").Value = "sd"
$ws.Range("J169").Value = "Statement-non-opinion"
$ws.Range("I172").Value = "sd"
$ws.Range("J172").Value = "Statement-non-opinion"
$ws.Range("I208").Value = "sd"
$ws.Range("J208").Value = "Statement-non-opinion"
$ws.Range("I219").Value = "sd"
$ws.Range("J219").Value = "Statement-non-opinion"
$ws.Range("I227").Value = "sd"
$ws.Range("J227").Value = "Statement-non-opinion"
$ws.Range("I230").Value = "aa"
$ws.Range("J230").Value = "Agree/Accept"
$ws.Range("I248").Value = "sd"
$ws.Range("J248").Value = "Statement-non-opinion"
$ws.Range("I252").Value = "sd"
$ws.Range("J252").Value = "Statement-non-opinion"
$ws.Range("I255").Value = "%"
$ws.Range("J255").Value = "Uninterpretable"
$ws.Range("I256").Value = "ba"
$ws.Range("J256").Value = "Appreciation"
$ws.Range("I258").Value = "sd"
$ws.Range("J258").Value = "Statement-non-opinion"
$ws.Range("I263").Value = "b"
$ws.Range("J263").Value = "Acknowledge (Backchannel)"
$ws.Range("I269").Value = "%"
$ws.Range("J269").Value = "Uninterpretable"
$ws.Range("I273").Value = "aa"
$ws.Range("J273").Value = "Agree/Accept"
$ws.Range("I275").Value = "sd"
$ws.Range("J275").Value = "Statement-non-opinion"
$ws.Range("I285").Value = "%"
$ws.Range("J285").Value = "Uninterpretable"
$ws.Range("I298").Value = "sd"
$ws.Range("J298").Value = "Statement-non-opinion"
$ws.Range("I300").Value = "sd"
$ws.Range("J300").Value = "Statement-non-opinion"
$ws.Range("I302").Value = "sd"
$ws.Range("J302").Value = "Statement-non-opinion"
$ws.Range("I304").Value = "sd"
$ws.Range("J304").Value = "Statement-non-opinion"
$ws.Range("I306").Value = "%"
$ws.Range("J306").Value = "Uninterpretable"
$ws.Range("I307").Value = "sv"
$ws.Range("J307").Value = "Statement-opinion"
$ws.Range("I313").Value = "sd"
$ws.Range("J313").Value = "Statement-non-opinion"
$ws.Range("I338").Value = "sd"
$ws.Range("J338").Value = "Statement-non-opinion"
$ws.Range("I355").Value = "sd"
$ws.Range("J355").Value = "Statement-non-opinion"
$ws.Range("I357").Value = "aa"
$ws.Range("J357").Value = "Agree/Accept"
$ws.Range("I358").Value = "aa"
$ws.Range("J358").Value = "Agree/Accept"
$ws.Range("I365").Value = "sd"
$ws.Range("J365").Value = "Statement-non-opinion"
$ws.Range("I367").Value = "sd"
$ws.Range("J367").Value = "Statement-non-opinion"
$ws.Range("I374").Value = "sd"
$ws.Range("J374").Value = "Statement-non-opinion"
$ws.Range("I376").Value = "sv"
$ws.Range("J376").Value = "Statement-opinion"
$ws.Range("I384").Value = "%"
$ws.Range("J384").Value = "Uninterpretable"
$ws.Range("I397").Value = "sd"
$ws.Range("J397").Value = "Statement-non-opinion"
$ws.Range("I403").Value = "sv"
$ws.Range("J403").Value = "Statement-opinion"
$ws.Range("I416").Value = "b"
$ws.Range("J416").Value = "Acknowledge (Backchannel)"
$ws.Range("I422").Value = "sv"
$ws.Range("J422").Value = "Statement-opinion"
$ws.Range("I423").Value = "b"
$ws.Range("J423").Value = "Acknowledge (Backchannel)"
$ws.Range("I426").Value = "%"
$ws.Range("J426").Value = "Uninterpretable"
$ws.Range("I431").Value = "sd"
$ws.Range("J431").Value = "Statement-non-opinion"
$ws.Range("I442").Value = "%"
$ws.Range("J442").Value = "Uninterpretable"
$ws.Range("I444").Value = "%"
$ws.Range("J444").Value = "Uninterpretable"
$ws.Range("I445").Value = "%"
$ws.Range("J445").Value = "Uninterpretable"
$ws.Range("I446").Value = "sv"
$ws.Range("J446").Value = "Statement-opinion"
$ws.Range("I448").Value = "aa"
$ws.Range("J448").Value = "Agree/Accept"
$ws.Range("I450").Value = "%"
$ws.Range("J450").Value = "Uninterpretable"
$ws.Range("I466").Value = "ba"
$ws.Range("J466").Value = "Appreciation"
$ws.Range("I473").Value = "sd"
$ws.Range("J473").Value = "Statement-non-opinion"
$ws.Range("I489").Value = "ba"
$ws.Range("J489").Value = "Appreciation"
$ws.Range("I500").Value = "aa"
$ws.Range("J500").Value = "Agree/Accept"
$ws.Range("I503").Value = "ba"
$ws.Range("J503").Value = "Appreciation"
$ws.Range("I518").Value = "sd"
$ws.Range("J518").Value = "Statement-non-opinion"
$ws.Range("I519").Value = "b"
$ws.Range("J519").Value = "Acknowledge (Backchannel)"
$ws.Range("I539").Value = "sv"
$ws.Range("J539").Value = "Statement-opinion"
$ws.Range("I540").Value = "sd"
$ws.Range("J540").Value = "Statement-non-opinion"
$ws.Range("I554").Value = "%"
$ws.Range("J554").Value = "Uninterpretable"
$ws.Range("I562").Value = "ba"
$ws.Range("J562").Value = "Appreciation"
$ws.Range("I567").Value = "sd"
$ws.Range("J567").Value = "Statement-non-opinion"
